$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update minhkhoi's "Tien" (money) value from 1000 to 1200
$ws.Range("B4").Value = 1200
